$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-11-07 07:03:31"

$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $range = $ws.Range("AA2:AA26")
    foreach ($cell in $range.Cells) {
        $cell.Value = $newTimestamp
    }
}
